$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 731 entirely (the "この広い宇宙に..." post), shifting rows 732-752 up to 731-751
$ws.Rows(731).Delete()
